$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: "accuracy" / temperatuur sensor -> temperatuur wasprogramma ---
$ws.Range("B2").Value = "temperatuur wasprogramma"
$ws.Range("C2").Value = "hoeveel de temperatuur mag afwijken van de ingestelden temperatuur"
$ws.Range("D2").Value = "maximaal 3 oC erboven en eronder vanaf ingestelde temperatuur"

# --- Row 3: "accuracy" / water sensor - update requirement text ---
$ws.Range("C3").Value = "hoeveel procent mag het vullen van de wasmachine trommel afwijken van de ingestelde hoeveelheid van het gekozen wasprogramma"
$ws.Range("D3").Value = "max 1% ofwel hij mag er 1% onder zitten of 1% overheen zitten."

# --- Row 4: becomes the former "usability / UI navigatie" row with reworded Beschrijving ---
$ws.Range("A4").Value = "usability"
$ws.Range("B4").Value = "UI navigatie"
$ws.Range("C4").Value = "het maximum aantal keer dat de gebruiker met de muis hoeft te klikken om in de goede volgorde het gewenste wasprogramma te starten"
$ws.Range("D4").Value = "3 - 4 muisklikken"
$ws.Range("E4").Value = "testen hoeveel acties er nodig zijn om een was te starten"
# the longer text now wraps across 4 lines, so the row grows from 40.5 to 51pt
$ws.Rows("4:4").RowHeight = 51

# --- Row 5: becomes the former "resource use" row with reworded Beschrijving ---
$ws.Range("A5").Value = "resource use"
$ws.Range("B5").Value = "ram gebruik besturing WM"
$ws.Range("C5").Value = "hoeveel ram mag het wasprogramma op de raspberry pi gebruiken"
$ws.Range("D5").Value = "max 120 MB"
$ws.Range("E5").Value = "door te testen en te meten hoeveel hij gebruikt tijdens uitvoeren"

# --- Row 6: becomes the former "performance / temperatuur weergeven" row with reworded Beschrijving ---
$ws.Range("A6").Value = "performance"
$ws.Range("B6").Value = "temperatuur weergeven"
$ws.Range("C6").Value = "Hoe snel en vaak moet de temperatuur weergaven geüpdate worden"
$ws.Range("D6").Value = "max 1 keer per seconde"
$ws.Range("E6").Value = "meten hoe veel tijd er tussen updates zit."

# --- Rows 7-9: old "vergrendelen WM deur" / "learnability" rows removed entirely ---
$ws.Range("A7:E9").Value = ""
$blankFormat = $ws.Range("G1")
$blankFormat.Copy()
$ws.Range("A7:E9").PasteSpecial(-4122)

# --- Final selection moves to D9 ---
$ws.Range("D9").Select() | Out-Null
